$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the commit diff (crypto price/volume refresh).
# Values that parse as plain numbers are written with a leading apostrophe
# (forcing text, matching the source inlineStr cells) and then have their
# formatting cleared so no stray quote-prefix style lingers on the cell.

$ws.Range('D2').Value = '63.100.82'
$ws.Range('E2').Value = '  +3.37%  '
$ws.Range('D3').Value = '3.033.26'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'593.21"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').Value = "'154.15"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +8.55%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '3.030.33'
$ws.Range('E8').Value = '  +1.95%  '
$ws.Range('D9').Value = "'0.516"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('D10').Value = "'6.86"
$ws.Range('D10').ClearFormats()
$ws.Range('E11').Value = '  +4.50%  '
$ws.Range('E12').Value = '  +2.08%  '
$ws.Range('E13').Value = '  +3.71%  '
$ws.Range('D14').Value = "'35.74"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +5.26%  '
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').Value = '3.537.90'
$ws.Range('E16').Value = '  +2.07%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '63.084.93'
$ws.Range('E17').Value = '  +3.35%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = "'7.08"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.28%  '
$ws.Range('D19').Value = '3.033.46'
$ws.Range('E19').Value = '  +2.13%  '
$ws.Range('D20').Value = "'453.23"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('D21').Value = "'14.27"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('E22').Value = '  +2.90%  '
$ws.Range('E23').Value = '  +3.40%  '
$ws.Range('D24').Value = "'83.10"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('E25').Value = '  +10.48%  '
$ws.Range('D26').Value = "'2.30"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.71%  '
$ws.Range('D27').Value = "'12.43"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.60%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = "'7.48"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +5.96%  '
$ws.Range('E30').Value = '  +11.23%  '
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').Value = "'0.999"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').Value = "'27.57"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.05%  '
$ws.Range('E34').Value = '  +2.37%  '
$ws.Range('D35').Value = '0.0₃0861'
$ws.Range('E35').Value = '  +6.60%  '
$ws.Range('E36').Value = '  +3.74%  '
$ws.Range('D37').Value = "'5.92"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.29%  '
$ws.Range('D38').Value = "'3.16"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +11.81%  '
$ws.Range('E39').Value = '  +8.91%  '
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('D41').Value = "'50.51"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range('D42').Value = "'9.12"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('E43').Value = '  +16.14%  '
$ws.Range('D44').Value = "'43.98"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +12.95%  '
$ws.Range('D45').Value = "'391.52"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.50%  '
$ws.Range('E46').Value = '  +3.83%  '
$ws.Range('D47').Value = '2.721.04'
$ws.Range('E47').Value = '  +1.73%  '
$ws.Range('D48').Value = "'133.39"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.41%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('E50').Value = '  +8.06%  '
$ws.Range('D51').Value = "'25.20"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +9.04%  '
